# Auto-generated script: apply 2025-01-21 crime data update
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 317
$ws.Range("K3").Value = 8183
$ws.Range("L3").Value = 313
$ws.Range("B4").Value = 1706
$ws.Range("K4").Value = 1730
$ws.Range("L4").Value = 79
$ws.Range("J6").Value = 11053
$ws.Range("K6").Value = 9136
$ws.Range("L6").Value = 399
$ws.Range("B7").Value = 23339
$ws.Range("J7").Value = 29319
$ws.Range("K7").Value = 27530
$ws.Range("L7").Value = 1131

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K2").Value = 240
$ws.Range("L6").Value = 11
$ws.Range("L7").Value = 43
$ws.Range("L8").Value = 62
$ws.Range("L11").Value = 21
$ws.Range("L14").Value = 4
$ws.Range("L19").Value = 43
$ws.Range("L20").Value = 26
$ws.Range("L23").Value = 12
$ws.Range("L25").Value = 7
$ws.Range("L27").Value = 9
$ws.Range("L29").Value = 63
$ws.Range("L31").Value = 14
$ws.Range("L33").Value = 45
$ws.Range("L36").Value = 21
$ws.Range("L37").Value = 39
$ws.Range("K43").Value = 227
$ws.Range("L44").Value = 6
$ws.Range("L49").Value = 9
$ws.Range("L52").Value = 23
$ws.Range("L53").Value = 15
$ws.Range("K57").Value = 112
$ws.Range("L58").Value = 1
$ws.Range("B63").Value = 410
$ws.Range("K63").Value = 76
$ws.Range("L63").Value = 6
$ws.Range("L64").Value = 12
$ws.Range("L65").Value = 18
$ws.Range("L67").Value = 31
$ws.Range("L70").Value = 6
$ws.Range("J76").Value = 413
$ws.Range("L76").Value = 19
$ws.Range("L77").Value = 9
$ws.Range("L78").Value = 12
$ws.Range("L79").Value = 29
$ws.Range("L83").Value = 27
$ws.Range("K85").Value = 1277
$ws.Range("L85").Value = 57
$ws.Range("L87").Value = 8
$ws.Range("L89").Value = 13
$ws.Range("K91").Value = 331
$ws.Range("L91").Value = 12
$ws.Range("L93").Value = 5
$ws.Range("L94").Value = 17
$ws.Range("L99").Value = 16
$ws.Range("B101").Value = 23339
$ws.Range("J101").Value = 29319
$ws.Range("K101").Value = 27530
$ws.Range("L101").Value = 1131

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("L6").Value = 3
$ws.Range("L7").Value = 4

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L2").Value = 10
$ws.Range("L3").Value = 12
$ws.Range("L7").Value = 43

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("L6").Value = 7
$ws.Range("L7").Value = 21

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("L2").Value = 4
$ws.Range("L7").Value = 13

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L3").Value = 23
$ws.Range("K4").Value = 63
$ws.Range("L6").Value = 12
$ws.Range("K7").Value = 1277
$ws.Range("L7").Value = 57

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("L2").Value = 7
$ws.Range("L7").Value = 23

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("L4").Value = 3
$ws.Range("L6").Value = 8
$ws.Range("L7").Value = 15

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L2").Value = 17
$ws.Range("L6").Value = 24
$ws.Range("L7").Value = 62

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("L2").Value = 9
$ws.Range("L7").Value = 27

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L3").Value = 16
$ws.Range("L6").Value = 21
$ws.Range("L7").Value = 45

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L2").Value = 14
$ws.Range("L7").Value = 39

$ws = $wb.Worksheets.Item('New City')
$ws.Range("L3").Value = 3
$ws.Range("L7").Value = 18

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("L3").Value = 5
$ws.Range("L7").Value = 16

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("L6").Value = 7
$ws.Range("L7").Value = 14

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L2").Value = 7
$ws.Range("L3").Value = 7
$ws.Range("L7").Value = 31

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("L6").Value = 4
$ws.Range("L7").Value = 9

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L2").Value = 21
$ws.Range("L4").Value = 3
$ws.Range("L6").Value = 22
$ws.Range("L7").Value = 63

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L2").Value = 17
$ws.Range("L3").Value = 9
$ws.Range("L6").Value = 16
$ws.Range("L7").Value = 43

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("L6").Value = 5
$ws.Range("L7").Value = 6

$ws = $wb.Worksheets.Item('River North')
$ws.Range("L3").Value = 2
$ws.Range("J6").Value = 213
$ws.Range("J7").Value = 413
$ws.Range("L7").Value = 19

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("L6").Value = 5
$ws.Range("L7").Value = 11

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("L3").Value = 5
$ws.Range("L6").Value = 4
$ws.Range("L7").Value = 12

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("L3").Value = 6
$ws.Range("L7").Value = 12

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("L2").Value = 7
$ws.Range("K4").Value = 21
$ws.Range("L6").Value = 3
$ws.Range("K7").Value = 331
$ws.Range("L7").Value = 12

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L3").Value = 9
$ws.Range("L7").Value = 29

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("L6").Value = 3
$ws.Range("L7").Value = 12

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("L6").Value = 8
$ws.Range("L7").Value = 26

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("L2").Value = 10
$ws.Range("L6").Value = 7
$ws.Range("L7").Value = 21

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("L6").Value = 3
$ws.Range("L7").Value = 5

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("L3").Value = 6
$ws.Range("L7").Value = 17

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("L3").Value = 3
$ws.Range("L7").Value = 7

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("K3").Value = 71
$ws.Range("K7").Value = 240

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("L2").Value = 2
$ws.Range("L3").Value = 3
$ws.Range("L7").Value = 6

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("L6").Value = 4
$ws.Range("L7").Value = 9

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("K3").Value = 22
$ws.Range("K7").Value = 112

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K6").Value = 82
$ws.Range("K7").Value = 227

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("L3").Value = 5
$ws.Range("L7").Value = 9

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range("L6").Value = 8
$ws.Range("L7").Value = 8

$ws = $wb.Worksheets.Item('Millenium Park')
$ws.Range("L1").Value = 2025
$ws.Range("L2").Value = 1
$ws.Range("L7").Value = 1
